# Added Backward extension option for real-time data
#
# The table of winter SMA YoY forecast values originally started in 1995
# (row 2). This adds 11 more years of backward-extended history
# (1984-1994) as new leading rows, pushing the pre-existing rows down by
# 11 (so the table now spans rows 2-42 instead of 2-31) and updates the
# sheet dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data that already existed in the sheet (previously rows 2-31), now
#     relocated to rows 13-42. Each tuple is (A,B,C,D,E).
$existingData = @(
    (35040, 1995, 2.234710814035812, 1996, 1.463127579670287),
    (35403, 1996, 1.595002781738275, 1997, 4.207635715208324),
    (35768, 1997, 2.499560583078497, 1998, 3.492506333467071),
    (36132, 1998, 2.812603855740181, 1999, 2.334197296693863),
    (36501, 1999, 1.188004848513446, 2000, 2.032004888754391),
    (36858, 2000, 3.277038745546235, 2001, 3.09884301635126),
    (37222, 2001, 0.7513248531724415, 2002, -0.408724114026926),
    (37581, 2002, 0.2537741062064169, 2003, 0.9118162660485263),
    (37938, 2003, -0.1535080579381121, 2004, 0.3435726964089891),
    (38302, 2004, 1.171834509066594, 2005, 0.8394840956263971),
    (38671, 2005, 1.120380359544382, 2006, 1.940699468213469),
    (39035, 2006, 2.691354324129258, 2007, 3.187301687590338),
    (39400, 2007, 2.652245539637632, 2008, 2.158031012958861),
    (39765, 2008, 1.327195601304898, 2009, -1.941693908020603),
    (40130, 2009, -4.803590807538871, 2010, 2.536922056245872),
    (40494, 2010, 3.776429555840499, 2011, 5.124900822223233),
    (40862, 2011, 3.167941427237042, 2012, 1.70423418303296),
    (41228, 2012, 1.072335020576287, 2013, 0.990934028412549),
    (41592, 2013, 0.5676944965793185, 2014, 1.859803271823757),
    (41957, 2014, 1.417171832295883, 2015, 0.05331272828721367),
    (42321, 2015, 1.475252114130599, 2016, 1.442973638880907),
    (42689, 2016, 1.71887541289224, 2017, 1.076548192761484),
    (43053, 2017, 2.581636142651922, 2018, 3.064375402422015),
    (43418, 2018, 1.471137749280693, 2019, 0.1137080120319656),
    (43783, 2019, 0.508332909595044, 2020, -0.1203207525434236),
    (44159, 2020, -5.494775307949129, 2021, 8.235743591092737),
    (44525, 2021, 3.149343082976164, 2022, 7.432336632701175),
    (44890, 2022, 1.995866057153428, 2023, 1.189587957345273),
    (45254, 2023, -0.1168430792840458, 2024, -0.1427298585871872),
    (45618, 2024, -0.1775688094211469, 2025, -0.1265568156813002)
)

# --- New backward-extension data (1984-1994 winter series), written into
#     the freshly-opened rows 2-12.
$newData = @(
    (31047, 1984, 2.833670241322217, 1985, 4.978977805976226),
    (31412, 1985, 2.740628897120945, 1986, 4.840042388885646),
    (31777, 1986, 2.269459987912947, 1987, 3.593781657196393),
    (32142, 1987, 1.253514454810789, 1988, 5.084502077712005),
    (32508, 1988, 3.509161092519553, 1989, 5.511076843601681),
    (32873, 1989, 3.898460078540933, 1990, 2.951715842334024),
    (33238, 1990, 5.356103277865332, 1991, 6.3181560832964),
    (33603, 1991, 5.955905607167122, 1992, -0.005898890116151634),
    (33969, 1992, 1.850401149566561, 1993, -0.4495646332120296),
    (34334, 1993, -0.9857661435315745, 1994, 2.795029892345036),
    (34699, 1994, 3.052254893522388, 1995, 3.383052772393214)
)

# 1) Re-write the pre-existing rows 11 positions further down (row 2 -> row 13,
#    row 31 -> row 42). The target rows already carry the correct formatting
#    (column A: date-formatted/bordered style, columns B:E: unstyled), so a
#    plain value write is all that's needed here.
for ($i = 0; $i -lt $existingData.Count; $i++) {
    $r = 13 + $i
    $row = $existingData[$i]
    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
    $ws.Cells.Item($r, 4).Value2 = $row[3]
    $ws.Cells.Item($r, 5).Value2 = $row[4]
}

# 2) Stamp rows 2-12 (which still hold stale copies of the old data at this
#    point) with the same cell formatting used throughout the table, sourced
#    from row 2's still-intact formatting.
$ws.Range("A2:E2").Copy()
$ws.Range("A13:E42").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# 3) Finally, fill rows 2-12 with the new backward-extension series.
for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = 2 + $i
    $row = $newData[$i]
    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
    $ws.Cells.Item($r, 4).Value2 = $row[3]
    $ws.Cells.Item($r, 5).Value2 = $row[4]
}
